# Gamelogic project removed; SLG building config (BB_Build) gains two new
# columns - "Icon" (short prefab/object name) and "ShowName" (display name,
# same text as the existing localized "Desc" column). The previous "Desc"
# column (G) shifts right to become column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 10

# 1) Shift the existing Desc column (G) out of the way into the new column I
#    first (walk bottom-up so we never clobber a value we still need to read).
for ($r = $lastRow; $r -ge 1; $r--) {
  $ws.Range("I$r").Value2 = $ws.Range("G$r").Value2
}

# 2) New header cells for the two inserted columns.
$ws.Range("G1").Value2 = "Icon"
$ws.Range("H1").Value2 = "ShowName"

# 3) Fill in the per-building data: Icon is the bare object name (the
#    "Prefabs/Object/" prefix stripped from column D), ShowName duplicates
#    the localized text that already lives in column I (old column G).
$icons = @(
  "Altar_1_1",
  "Arena_1_1",
  "Camp_1_1",
  "GoldMine_1_1",
  "Item_hourse_1_1",
  "League_1_1",
  "MagicHourse_1_1",
  "Tower_1_1",
  "Town_1_1"
)

for ($i = 0; $i -lt $icons.Length; $i++) {
  $r = $i + 2
  $ws.Range("G$r").Value2 = $icons[$i]
  $ws.Range("H$r").Value2 = $ws.Range("I$r").Value2
  # Match the text-formatted style already used by the rest of the row.
  $ws.Range("H$r").NumberFormat = "@"
  $ws.Range("I$r").NumberFormat = "@"
}

# 4) Column widths: G:I all become a uniform "11 chars" wide.
$ws.Range("G1:I1").EntireColumn.ColumnWidth = 10.285714285714286

# 5) Selection moves to H10 (where the last new value was entered).
$ws.Range("H10").Select() | Out-Null
